# Add the "Día 2" narrative to the end of the "Guía del proyecto" document.
#
# The document currently ends with a paragraph describing "Día 1" work; we
# append two new paragraphs ("Día 2:" and the description of what was done
# that day), and make sure the "_GoBack" bookmark (which Word keeps at the
# location of the most recent edit) ends up wrapping the end of the new
# last paragraph, matching what real Word does when you type new text at
# the end of a document.

$d = $word.ActiveDocument

# The last paragraph currently in the document ("En la clase de Taller,
# diseñamos ..."). We'll append our new paragraphs right after it.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)

# Collapse to just before that paragraph's end-of-paragraph mark and insert
# a new (empty) paragraph after it.
$insertionPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)
$insertionPoint.InsertParagraphAfter()

# Fill in the new "Día 2:" paragraph.
$day2Para = $d.Paragraphs($d.Paragraphs.Count)
$day2Para.Range.InsertAfter("Día 2:")

# Add another new paragraph after that one for the Día 2 description.
$day2Para = $d.Paragraphs($d.Paragraphs.Count)
$insertionPoint2 = $d.Range($day2Para.Range.End - 1, $day2Para.Range.End - 1)
$insertionPoint2.InsertParagraphAfter()

$descPara = $d.Paragraphs($d.Paragraphs.Count)
# Append a placeholder character temporarily so we can anchor the "_GoBack"
# bookmark right after the real text (but before the placeholder), then
# remove the placeholder. Placing the bookmark directly at the very end of
# the document's content can otherwise snap it back to the start of the
# document, so this keeps it anchored correctly at the end of the new text.
$descPara.Range.InsertAfter("En la clase de Taller, seguimos agregando precios de los objetos que escogimos y de ya tenemos los precios de los cables que vamos a utilizar en la siguiente clase.X")

$descPara = $d.Paragraphs($d.Paragraphs.Count)
$bookmarkPos = $d.Range($descPara.Range.End - 2, $descPara.Range.End - 2)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bookmarkPos)

$placeholder = $d.Range($descPara.Range.End - 2, $descPara.Range.End - 1)
$placeholder.Delete()
